# Benchmark update: 2026-01-15 06:43:29 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear stale benchmark figures (keep cell style, drop text) ---
$clearAddresses = @(
    "D3","E3","F3","I3",
    "D4","E4","F4","I4",
    "D5","E5","F5","I5",
    "D6","E6","I6",
    "D8","E8","F8","I8",
    "D9","E9","F9","I9",
    "D10","E10","F10","I10",
    "D11","E11","I11",
    "D12",
    "D13","E13","F13","I13",
    "D14","E14","F14",
    "K20",
    "K23",
    "K24",
    "K25"
)

foreach ($addr in $clearAddresses) {
    $ws.Range($addr).ClearContents()
}

# --- Updated benchmark figures ---
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("K14").Value = "2.000 TL - 24.000 TL"
$ws.Range("K15").Value = " Asgari Tutar:  Azami Tutar: "
$ws.Range("K17").Value = " Asgari Tutar:  Azami Tutar: "
$ws.Range("K21").Value = " Asgari Tutar:  Azami Tutar: "
$ws.Range("K22").Value = " Asgari Tutar:  Azami Tutar: "
